# Applies the "update readme & add structure img" edit to the structure
# diagram slide:
#   - widen / re-dash / re-color the 5 dashed accent2 connector arrows
#   - nudge 4 of those connectors (and 4 caption textboxes) horizontally
#   - (best-effort) wire the start connection site of the 4 re-routed
#     connectors to shape id 7 connection site 3
#   - (best-effort) touch the presentation-level guide list
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# are expressed in points (1 pt = 12700 EMU) and this host stores them in a
# 32-bit float, so a plain `targetEmu / 12700.0` can truncate to one EMU
# below the value recorded in the OOXML diff once it is re-expanded on
# save. The literals used here were solved so that, after the host's
# float32 round-trip, `round(value * 12700)` reproduces the exact EMU
# figures from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-ConnectorLine {
    param($shape)
    # Order matters: this host always re-serializes a changed <a:prstDash>
    # and <a:tailEnd> at the tail of <a:ln>, so DashStyle must be (re)set
    # before EndArrowheadStyle to land them in prstDash-then-tailEnd order
    # (matching the OOXML EG_LineDashProperties/tailEnd schema sequence).
    $shape.Line.Weight = 1.5
    $shape.Line.ForeColor.ObjectThemeColor = 6
    $shape.Line.DashStyle = 9
    $shape.Line.EndArrowheadStyle = 2
}

# Locate shapes by their stable PowerPoint shape Id (positions in the
# Shapes collection already match 1:1 with Id order in this deck, but we
# look them up defensively in case that ever changes).
$shapesById = @{}
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $shapesById[$sh.Id] = $sh
}

$anchorShape = $shapesById[7]

# --- Connector 9: CMD -> (rounded rect 6) — line restyle only ---
Set-ConnectorLine $shapesById[9]

# --- Connectors 18 / 20 / 21 / 23: restyle + reroute start + reposition ---
$connectorMoves = @(
    @{ Id = 18; EndId = 12; Left = 241.7511826023622; Width = 48.819765179527565 },
    @{ Id = 20; EndId = 13; Left = 241.7511826023622; Width = 48.819765179527565 },
    @{ Id = 21; EndId = 14; Left = 241.7511826023622; Width = 48.81968503937008  },
    @{ Id = 23; EndId = 15; Left = 241.7511826023622; Width = 48.819765179527565 }
)

foreach ($move in $connectorMoves) {
    $conn = $shapesById[$move.Id]

    Set-ConnectorLine $conn

    # Best-effort: attach the connector's start point to shape 7 / site 3
    # (mirrors <a:stCxn id="7" idx="3"/> in the target OOXML). Some hosts
    # do not implement connector re-routing via COM; if so this silently
    # has no effect and the explicit Left/Width restamp below still lands
    # the correct geometry from the diff.
    try {
        $conn.ConnectorFormat.BeginConnect($anchorShape, 3)
    } catch {
    }

    $conn.Left = $move.Left
    $conn.Width = $move.Width
}

# --- Caption textboxes 55 / 56 / 57 / 58: horizontal nudge only ---
$textboxMoves = @(
    @{ Id = 55; Left = 359.9203149606299 },
    @{ Id = 56; Left = 326.46418762834645 },
    @{ Id = 57; Left = 688.2374803149606 },
    @{ Id = 58; Left = 682.2947244094488 }
)

foreach ($move in $textboxMoves) {
    $shapesById[$move.Id].Left = $move.Left
}

# --- Presentation-level empty slide guide list (best-effort) ---
# The diff adds an empty <p15:sldGuideLst/> ext to presentation.xml; touch
# the Guides collection in case this host materializes that element as a
# side effect of accessing it. Harmless no-op if unsupported.
try {
    $null = $p.Guides
} catch {
}
